$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# B14 ("Case Sensitive" value) becomes the same text value currently held by
# B7 ("Experimental" = "false", stored as text/shared-string, not boolean).
# Copy/paste (rather than typing) keeps it text instead of letting Excel's
# TRUE/FALSE auto-typing turn it into a Boolean cell.
$ws.Range("B7").Copy($ws.Range("B14"))

# Flip "Experimental" (and now also "Case Sensitive") from "false" to "true".
# Build the literal through a formula and paste its *value* so Excel stores
# plain text "true" instead of auto-converting it to a Boolean.
$ws.Range("D1").Formula = "=""tr""&""ue"""
$ws.Range("D1").Copy()
$ws.Range("B7").PasteSpecial(-4163)
$ws.Range("B14").PasteSpecial(-4163)
$ws.Range("D1").Clear()

# Update the Date metadata value.
$ws.Range("B8").Value = "2024-02-19T18:37:26-06:00"
